$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on cells whose new values look numeric, to preserve
# the original text representation (matching trailing zeros, etc.)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated values
$ws.Range("D2").Value = '42.930.45'
$ws.Range("E2").Value = '  +0.53%  '
$ws.Range("D3").Value = '2.361.70'
$ws.Range("E3").Value = '  +2.23%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = '302.10'
$ws.Range("E5").Value = '  +0.19%  '
$ws.Range("D6").Value = '95.77'
$ws.Range("E6").Value = '  +0.43%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").Value = '0.501'
$ws.Range("E8").Value = '  -0.39%  '
$ws.Range("D9").Value = '0.488'
$ws.Range("E9").Value = '  -0.53%  '
$ws.Range("D10").Value = '34.09'
$ws.Range("E10").Value = '  -0.15%  '
$ws.Range("D11").Value = '0.0783'
$ws.Range("E11").Value = '  +0.18%  '
$ws.Range("D12").Value = '0.124'
$ws.Range("E12").Value = '  +3.32%  '
$ws.Range("D13").Value = '18.31'
$ws.Range("E13").Value = '  -2.86%  '
$ws.Range("D14").Value = '6.72'
$ws.Range("E14").Value = '  +0.01%  '
$ws.Range("D15").Value = '2.731.06'
$ws.Range("E15").Value = '  +2.25%  '
$ws.Range("D16").Value = '2.357.94'
$ws.Range("E16").Value = '  +1.84%  '
$ws.Range("D17").Value = '0.791'
$ws.Range("E17").Value = '  +0.42%  '
$ws.Range("D18").Value = '42.900.79'
$ws.Range("E18").Value = '  +0.59%  '
$ws.Range("B19").Value = 'Uniswap'
$ws.Range("C19").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D19").Value = '6.24'
$ws.Range("E19").Value = '  +2.01%  '
$ws.Range("B20").Value = 'InternetComputer(DFINITY)'
$ws.Range("C20").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D20").Value = '11.81'
$ws.Range("E20").Value = '  -2.62%  '
$ws.Range("D21").Value = '0.0₃0882'
$ws.Range("D22").Value = '67.85'
$ws.Range("E22").Value = '  +0.22%  '
$ws.Range("D23").Value = '234.89'
$ws.Range("E23").Value = '  -0.08%  '
$ws.Range("D24").Value = '2.15'
$ws.Range("E24").Value = '  -5.00%  '
$ws.Range("E25").Value = '  +0.70%  '
$ws.Range("E26").Value = '  -0.06%  '
$ws.Range("D27").Value = '24.45'
$ws.Range("E27").Value = '  +0.76%  '
$ws.Range("E28").Value = '  +0.88%  '
$ws.Range("E29").Value = '  +1.97%  '
$ws.Range("D30").Value = '31.84'
$ws.Range("E30").Value = '  -0.80%  '
$ws.Range("E31").Value = '  -0.02%  '
$ws.Range("E32").Value = '  +0.31%  '
$ws.Range("D33").Value = '17.26'
$ws.Range("E33").Value = '  -1.90%  '
$ws.Range("D34").Value = '0.0710'
$ws.Range("E34").Value = '  +1.97%  '
$ws.Range("E35").Value = '  +3.86%  '
$ws.Range("D36").Value = '1.84'
$ws.Range("E36").Value = '  +3.40%  '
$ws.Range("D37").Value = '4.33'
$ws.Range("E37").Value = '  -2.92%  '
$ws.Range("D38").Value = '125.51'
$ws.Range("E38").Value = '  -24.63%  '
$ws.Range("E39").Value = '  -1.89%  '
$ws.Range("E40").Value = '  +3.78%  '
$ws.Range("D41").Value = '0.107'
$ws.Range("E41").Value = '  -0.85%  '
$ws.Range("D42").Value = '21.20'
$ws.Range("E42").Value = '  +0.32%  '
$ws.Range("D43").Value = '1.929.03'
$ws.Range("E43").Value = '  +0.24%  '
$ws.Range("D44").Value = '0.0277'
$ws.Range("E44").Value = '  -0.34%  '
$ws.Range("E45").Value = '  +2.34%  '
$ws.Range("E46").Value = '  -0.53%  '
$ws.Range("D47").Value = '9.16'
$ws.Range("E47").Value = '  -8.67%  '
$ws.Range("D48").Value = '2.589.97'
$ws.Range("E48").Value = '  +1.97%  '
$ws.Range("E49").Value = '  +1.86%  '
$ws.Range("E50").Value = '  +1.40%  '
$ws.Range("D51").Value = '51.47'
$ws.Range("E51").Value = '  -3.47%  '
